# Applies the "lecon 4 et debut 5" update to the horaire worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (cours #4): was just "Express", becomes a generator lesson with a new exercise reference.
# (Set D5 before C5 so the shared-string table gets the two new strings in the same order
# as the authoritative workbook: exercise text first, then the subject text.)
$ws.Range("D5").Value = "[Exercice 4 - Express avec modèle](exercice4_express_avec_modele.md)"
$ws.Range("C5").Value = "[Express - Générateur](generateur_express.md)"

# Row 6 (cours #5): was just "MongoDB", becomes async JS + MongoDB with its own exercise reference.
$ws.Range("C6").Value = "[JavaScript asynchrone](javascript_async.md) <br/> [MongoDB](mongodb.md)"
$ws.Range("D6").Value = "[Exercice 5 - MongoDB](exercice5_mongodb.md)"

# Update the active selection to reflect where the author left off editing.
$ws.Range("E5").Select()
